$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vt = [char]11

$cell = $t.Cell(1,1)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "95 x 25" + $vt + "  2    5" + $vt + "  ----" + $vt + "9|    |" + $vt + "5|    |"

$cell = $t.Cell(1,2)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "54 x 37" + $vt + "  3    7" + $vt + "  ----" + $vt + "5|    |" + $vt + "4|    |"

$cell = $t.Cell(1,3)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "29 x 33" + $vt + "  3    3" + $vt + "  ----" + $vt + "2|    |" + $vt + "9|    |"

$cell = $t.Cell(2,1)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "28 x 55" + $vt + "  5    5" + $vt + "  ----" + $vt + "2|    |" + $vt + "8|    |"

$cell = $t.Cell(2,2)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "19 x 67" + $vt + "  6    7" + $vt + "  ----" + $vt + "1|    |" + $vt + "9|    |"

$cell = $t.Cell(2,3)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "73 x 88" + $vt + "  8    8" + $vt + "  ----" + $vt + "7|    |" + $vt + "3|    |"

$cell = $t.Cell(3,1)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "41 x 89" + $vt + "  8    9" + $vt + "  ----" + $vt + "4|    |" + $vt + "1|    |"

$cell = $t.Cell(3,2)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "27 x 68" + $vt + "  6    8" + $vt + "  ----" + $vt + "2|    |" + $vt + "7|    |"

$cell = $t.Cell(3,3)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "75 x 87" + $vt + "  8    7" + $vt + "  ----" + $vt + "7|    |" + $vt + "5|    |"

$cell = $t.Cell(4,1)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "86 x 54" + $vt + "  5    4" + $vt + "  ----" + $vt + "8|    |" + $vt + "6|    |"

$cell = $t.Cell(4,2)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "60 x 19" + $vt + "  1    9" + $vt + "  ----" + $vt + "6|    |" + $vt + "0|    |"

$cell = $t.Cell(4,3)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "45 x 86" + $vt + "  8    6" + $vt + "  ----" + $vt + "4|    |" + $vt + "5|    |"

$cell = $t.Cell(5,1)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "69 x 89" + $vt + "  8    9" + $vt + "  ----" + $vt + "6|    |" + $vt + "9|    |"

$cell = $t.Cell(5,2)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "22 x 12" + $vt + "  1    2" + $vt + "  ----" + $vt + "2|    |" + $vt + "2|    |"

$cell = $t.Cell(5,3)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "87 x 35" + $vt + "  3    5" + $vt + "  ----" + $vt + "8|    |" + $vt + "7|    |"
